$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (week number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/4/2024  Through  11/10/2024"

# --- Cells changing from blank-placeholder text ("0"/"***.*") to real numbers ---
# Setting NumberFormat first makes the runtime reuse the existing numeric
# style (15 for counts, 14 for percent-change) instead of minting a new one.
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("D16").Value = 2
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E16").Value = -50
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 1
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = -50
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 1
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H22").Value = -100
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 1
$ws.Range("F23").NumberFormat = "#,##0"
$ws.Range("F23").Value = 1

# --- Cells changing from real numbers to placeholder text ("0"/"***.*") ---
# Copying from an existing placeholder cell elsewhere on the sheet reproduces
# the exact shared-string + style pairing Excel itself uses for these cells.
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))

# --- Plain numeric value updates ---
$ws.Range("L15").Value = -30
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = 8.510638297872
$ws.Range("M16").Value = -12.068965517241
$ws.Range("N16").Value = -78.925619834710
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -60
$ws.Range("I17").Value = 136
$ws.Range("J17").Value = 112
$ws.Range("K17").Value = 21.428571428571
$ws.Range("L17").Value = 46.236559139784
$ws.Range("M17").Value = 76.623376623376
$ws.Range("N17").Value = -29.166666666666
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = 25.641025641025
$ws.Range("L18").Value = -5.769230769230
$ws.Range("M18").Value = -3.921568627450
$ws.Range("N18").Value = -86.756756756756
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 7
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = -61.111111111111
$ws.Range("I19").Value = 137
$ws.Range("J19").Value = 168
$ws.Range("K19").Value = -18.452380952381
$ws.Range("L19").Value = -5.517241379310
$ws.Range("M19").Value = 33.009708737864
$ws.Range("N19").Value = -23.033707865168
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -66.666666666666
$ws.Range("L20").Value = -14.285714285714
$ws.Range("N20").Value = -83.448275862069
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 21
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = -52.272727272727
$ws.Range("I21").Value = 428
$ws.Range("J21").Value = 392
$ws.Range("K21").Value = 9.183673469387
$ws.Range("L21").Value = 5.940594059405
$ws.Range("M21").Value = 34.591194968553
$ws.Range("N21").Value = -66.795965865011
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 116.666666666667
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 21
$ws.Range("K23").Value = 90.909090909090
$ws.Range("L23").Value = 16.666666666666
$ws.Range("M23").Value = 10.526315789473
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = 5.555555555555
$ws.Range("I24").Value = 441
$ws.Range("J24").Value = 419
$ws.Range("K24").Value = 5.250596658711
$ws.Range("L24").Value = -5.769230769230
$ws.Range("M24").Value = 67.680608365019
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 76.923076923076
$ws.Range("I25").Value = 227
$ws.Range("J25").Value = 217
$ws.Range("K25").Value = 4.608294930875
$ws.Range("L25").Value = 11.822660098522
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = -29.411764705882
$ws.Range("I26").Value = 229
$ws.Range("J26").Value = 195
$ws.Range("K26").Value = 17.435897435897
$ws.Range("L26").Value = 20.526315789473
$ws.Range("M26").Value = -24.172185430463
$ws.Range("L27").Value = -23.076923076923
$ws.Range("L28").Value = -19.047619047619
$ws.Range("L29").Value = -12.5
$ws.Range("L30").Value = -16.666666666666
